# Daily automated EPEX Spot price update.
# - "Prix Spot": append a new day column (CW) with header "22-sep" and its 24 hourly prices.
# - "Gaz" / "CO2": append two new daily rows (2025-09-20, 2025-09-21).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": new column CW (column 101) after existing CV (column 100)
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

$lastCol = 100          # CV
$newCol  = $lastCol + 1 # CW

# Clone the header cell's formatting (bold, centered, bordered) from CV1 onto
# CW1, then set its text - mirrors the existing header style exactly.
$wsSpot.Cells.Item(1, $lastCol).Copy()
$wsSpot.Cells.Item(1, $newCol).PasteSpecial(-4122)
$wsSpot.Cells.Item(1, $newCol).Value = "22-sep"

# Hourly values for 22-sep, row 2 = "00 - 01" ... row 25 = "23 - 24"
$spotValues = @(
    @(2, 15.57),
    @(3, 14.08),
    @(4, 21.44),
    @(5, 11.33),
    @(6, 4.14),
    @(7, 1.5),
    @(8, 10.66),
    @(9, 41.9),
    @(10, 62.03),
    @(11, 29.35),
    @(12, 2),
    @(13, 0),
    @(14, 0),
    @(15, -0.01),
    @(16, -0.01),
    @(17, -0.01),
    @(18, -0.01),
    @(19, -0.01),
    @(20, 1.54),
    @(21, 39.66),
    @(22, 17.28),
    @(23, 9.22),
    @(24, 6.7),
    @(25, 26.47)
)

foreach ($pair in $spotValues) {
    $row = $pair[0]
    $val = $pair[1]
    $wsSpot.Cells.Item($row, $newCol).Value = $val
}

# ---------------------------------------------------------------------------
# Sheets "Gaz" and "CO2": two new trailing rows, 98 and 99.
# The date column stores plain text (e.g. "2025-09-20"), matching the rest of
# column A. A leading apostrophe forces text entry instead of Excel's
# automatic date conversion.
# ---------------------------------------------------------------------------
$newRows = @(
    @(98, "2025-09-20"),
    @(99, "2025-09-21")
)

$wsGaz = $wb.Worksheets.Item("Gaz")
$gazValues = @{ 98 = 31.75; 99 = 31.75 }
foreach ($pair in $newRows) {
    $row  = $pair[0]
    $date = $pair[1]
    $wsGaz.Cells.Item($row, 1).Value = "'" + $date
    $wsGaz.Cells.Item($row, 2).Value = $gazValues[$row]
}

$wsCo2 = $wb.Worksheets.Item("CO2")
$co2Values = @{ 98 = 76.63; 99 = 76.63 }
foreach ($pair in $newRows) {
    $row  = $pair[0]
    $date = $pair[1]
    $wsCo2.Cells.Item($row, 1).Value = "'" + $date
    $wsCo2.Cells.Item($row, 2).Value = $co2Values[$row]
}
